$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "فاطمه شادي عايش محمد "
$ws.Range("B2").Value = "fatmashady18@gmail.com"
$ws.Range("C2").Value = "https://github.com/Fatmashady/Security-Task.git"

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/Fatmashady/Security-Task.git")

$ws.Range("C4").Select()
